$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("acharyan_captions")
$ws2 = $wb.Worksheets.Item("Founders_Early_Acharyas")
$ws3 = $wb.Worksheets.Item("Banner_Text")

# --- Sheet1 (acharyan_captions): normalize "Mahā Deśikan" -> "Mahādeśikan" ---
$ws1.Range("B3").Value = "Sri Dviteeya (Vātsya Vedānta Rāmānuja) Brahmatantra Swatantra Mahādeśikan (1386 – 1394)"
$ws1.Range("B4").Value = "Sri Tṛtīya (Srinivāsa) Brahmatantra Swatantra Mahādeśikan (1394 – 1406)"
$ws1.Range("B5").Value = "Sri Parakāla Brahmatantra Swatantra Mahādeśikan (1406 – 1424)"
$ws1.Range("B6").Value = "Sri Vedānta Rāmānuja Brahmatantra Swatantra Mahādeśikan (1424 – 1440)"
$ws1.Range("B7").Value = "Sri Srinivāsa Brahmatantra Swatantra Mahādeśikan (1440 – 1460)"
$ws1.Range("B8").Value = "Sri Nārāyaṇa Yogindra Brahmatantra Swatantra Mahādeśikan (1460 – 1482)"
$ws1.Range("B9").Value = "Sri Raṅgarāja Yogindra Brahmatantra Swatantra Mahādeśikan (1482 – 1498)"
$ws1.Range("B10").Value = "Sri Chaturtha Brahmatantra Swatantra Mahādeśikan (1498 – 1517)"
$ws1.Range("B11").Value = "Sri Yatirāja Brahmatantra Swatantra Mahādeśikan (1517 – 1535)"
$ws1.Range("B12").Value = "Sri Varada Brahmatantra Swatantra Mahādeśikan (1535 – 1552)"
$ws1.Range("B13").Value = "Sri Parāṅkuśa Brahmatantra Swatantra Mahādeśikan (1552 – 1567)"
$ws1.Range("B14").Value = "Sri Kavitārkika Siṃha Brahmatantra Swatantra Mahādeśikan (1567 – 1583)"
$ws1.Range("B15").Value = "Sri Vedānta Yathivarya Brahmatantra Swatantra Mahādeśikan (1583 – 1607)"
$ws1.Range("B16").Value = "Sri Jñānābdi Brahmatantra Swatantra Mahādeśikan (1607 – 1618)"
$ws1.Range("B17").Value = "Sri Vīrarāghava Yogindra Brahmatantra Swatantra Mahādeśikan (1618 – 1640)"
$ws1.Range("B18").Value = "Sri Varada Vedānta Brahmatantra Swatantra Mahādeśikan (1640 – 1652)"
$ws1.Range("B19").Value = "Sri Varāha Brahmatantra Swatantra Mahādeśikan (1652 – 1663)"
$ws1.Range("B20").Value = "Sri Vedānta Lakṣmaṇa Brahmatantra Swatantra Mahādeśikan (1663 – 1673)"
$ws1.Range("B21").Value = "Sri Varada Vedānta Yogīndra Brahmatantra Swatantra Mahādeśikan (1673 – 1677)"
$ws1.Range("B22").Value = "Sri Maha Parakāla Brahmatantra Swatantra Mahādeśikan (1676 – 1738)"
$ws1.Range("B23").Value = "Sri Srinivāsa Brahmatantra Swatantra Parakāla Mahādeśikan (1738 – 1751)"
$ws1.Range("B24").Value = "Sri Vedānta Brahmatantra Swatantra Parakāla Mahādeśikan (1750 – 1770)"
$ws1.Range("B25").Value = "Srimat Abhinava Srinivāsa Brahmatantra Swatantra Parakāla Mahādeśikan (1770 – 1772)"
$ws1.Range("B27").Value = "Sri Ghantāvatāra Brahmatantra Swatantra Parakāla Mahādeśikan (1810 – 1829)"
$ws1.Range("B28").Value = "Sri Vedānta Brahmatantra Swatantra Parakāla Mahādeśikan (1829 – 1836)"
$ws1.Range("B29").Value = "Sri Srinivāsa Brahmatantra Swatantra Parakāla Mahādeśikan (1836 – 1861)"
$ws1.Range("B30").Value = "Sri Srinivāsa Deśikendra Brahmatantra Swatantra Parakāla Mahādeśikan (1861 – 1873)"
$ws1.Range("B31").Value = "Sri Raṅganātha Brahmatantra Swatantra Parakāla Mahādeśikan (1873 – 1885)"
$ws1.Range("B32").Value = "Sri Kṛṣṇa Brahmatantra Swatantra Parakāla Mahādeśikan (1885 – 1915)"
$ws1.Range("B33").Value = "Sri Vāgīśa Brahmatantra Swatantra Parakāla Mahādeśikan (1915 – 1925)"
$ws1.Range("B34").Value = "Srimat Abhinava Raṅganātha Brahmatantra Swatantra Parakāla Mahādeśikan (1925 – 1967)"
$ws1.Range("B35").Value = "Srimat Abhinava Srinivāsa Brahmatantra Swatantra Parakāla Mahādeśikan (1967 – 1972)"
$ws1.Range("B36").Value = "Srimat Abhinava Rāmānuja Brahmatantra Swatantra Parakāla Mahādeśikan (1972 – 1992)"
$ws1.Range("B37").Value = "Srimat Abhinava Vāgīśa Brahmatantra Swatantra Parakāla Mahādeśikan (1992 – present)"

# --- Sheet2 (Founders_Early_Acharyas): fix "Nigamānta Mahā Deśikan" spelling ---
$ws2.Range("C19").Value = "Nigamānta Mahādeśikan"

# --- Sheet3 (Banner_Text): rename M1 caption, clear the " Maṭham " column D cells ---
$ws3.Range("B2").Value = "Sri Brahmatantra Swatantra Parakala Swāmy Guru Parampara"
$ws3.Range("D2").Value = ""
$ws3.Range("D4").Value = ""
$ws3.Range("D6").Value = ""

# --- Recreate the view/selection state from the target workbook ---
$ws3.Activate()
$ws3.Range("D:D").Select()

$ws2.Activate()
$ws2.Range("C19").Select()

$ws1.Activate()
$ws1.Range("B5").Select()